# Applies two changes from the source diff:
#   1. The footer "datetimeFigureOut" date field on the slide master and on
#      every slide layout changes its displayed text from "10/01/2023" to
#      "2023-10-05".
#   2. The "Oval 20" shape (nested inside the top-level "Group 2" group on
#      slide 1) moves from local offset (2774597, 727587) EMU to
#      (2372410, 875230) EMU.

$p = $ppt.ActivePresentation

# --- 1. Update the date placeholder text everywhere it appears ---------
$newDate = "2023-10-05"
$ppPlaceholderDate = 16

function Update-DatePlaceholders($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master footer date placeholder.
Update-DatePlaceholders $p.SlideMaster.Shapes

# Every slide layout's footer date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholders $layouts.Item($li).Shapes
}

# --- 2. Reposition the "Oval 20" shape nested inside "Group 2" ---------
$slide = $p.Slides.Item(1)
$group = $slide.Shapes.Item(1)

for ($gi = 1; $gi -le $group.GroupItems.Count; $gi++) {
    $item = $group.GroupItems.Item($gi)
    if ($item.Name -eq "Oval 20") {
        $item.Left = 2372410 / 12700.0
        $item.Top = 875230 / 12700.0
    }
}
